# Update scripts with new TPM values (NATMI LR-pairs recomputation)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 3.027114666666666
$ws.Range("H2").Value = 9.081344
$ws.Range("I2").Value = 0.207506525262911
$ws.Range("J2").Value = 0.207506525262911
$ws.Range("M2").Value = 0.073876
$ws.Range("Q2").Value = 0.2236311231146667
$ws.Range("R2").Value = 2.012680108032
$ws.Range("S2").Value = 0.207506525262911
$ws.Range("T2").Value = 0.207506525262911

# Row 3
$ws.Range("H3").Value = 7.555711000000001
$ws.Range("I3").Value = 0.1726461783080517
$ws.Range("J3").Value = 0.1726461783080516
$ws.Range("M3").Value = 0.073876
$ws.Range("Q3").Value = 0.1860619019453333
$ws.Range("S3").Value = 0.1726461783080517
$ws.Range("T3").Value = 0.1726461783080516

# Row 4
$ws.Range("G4").Value = 4.235286666666666
$ws.Range("H4").Value = 12.70586
$ws.Range("I4").Value = 0.2903258437382188
$ws.Range("J4").Value = 0.2903258437382187
$ws.Range("M4").Value = 0.073876
$ws.Range("Q4").Value = 0.3128860377866666
$ws.Range("R4").Value = 2.81597434008
$ws.Range("S4").Value = 0.2903258437382188
$ws.Range("T4").Value = 0.2903258437382187

# Row 5
$ws.Range("G5").Value = 1.937427333333333
$ws.Range("H5").Value = 5.812282
$ws.Range("I5").Value = 0.1328092451588843
$ws.Range("J5").Value = 0.1328092451588843
$ws.Range("M5").Value = 0.073876
$ws.Range("Q5").Value = 0.1431293816773333
$ws.Range("R5").Value = 1.288164435096
$ws.Range("S5").Value = 0.1328092451588843
$ws.Range("T5").Value = 0.1328092451588843

# Row 6
$ws.Range("G6").Value = 2.869646666666667
$ws.Range("H6").Value = 8.60894
$ws.Range("I6").Value = 0.1967122075319342
$ws.Range("J6").Value = 0.1967122075319342
$ws.Range("M6").Value = 0.073876
$ws.Range("Q6").Value = 0.2119980171466667
$ws.Range("R6").Value = 1.90798215432
$ws.Range("S6").Value = 0.1967122075319342
$ws.Range("T6").Value = 0.1967122075319342
